# "rolled out new secure model"
#
# The template's data block (header row + sample row + mandatory-field
# marker row) previously started on row 2, leaving row 1 blank. The new
# version removes that leading blank row so the header row becomes row 1,
# and corrects the "YES NO" label (shared string) to "YES_NO".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank leading row so data starts at row 1 (shifts
# rows 2-4 up to rows 1-3, preserving all values/styles).
$ws.Rows("1:1").Delete()

# The ACTIVATION_STATUS header cell (now H1) held the shared string
# "YES NO" - correct it to "YES_NO".
$ws.Range("H1").Value = "YES_NO"
